# Adds new "Subrogation Loss" (subrogation expense / cause of loss) related
# keys/XPaths to the Navigation sheet, right after the existing
# "SubrogateAmtSubmitbtn" entry (row 83), extending the used range to F88.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows, columns D (Key) and E (XPath/Value), matching the pattern of the
# preceding Subrogation rows (80-83).
$newRows = @(
    @{ Row = 84; D = "SubrogateLossRadioBox"; E = "/html/body/div[2]/div/form/div/div[2]/div[1]/div/div[1]/div/div[2]/div/div[1]/div/div/div[3]/div/div/div/div[6]/div/div/div[2]/div[1]/div/table/tbody/tr[4]/td[1]/div/input" },
    @{ Row = 85; D = "SubrogateLossSubmit";   E = "/html/body/div[2]/div/form/div/div[2]/div[1]/div/div[1]/div/div[2]/div/div[1]/div/div/div[3]/div/div/div/div[6]/div/div/div[2]/div[2]/button[1]" },
    @{ Row = 86; D = "SubrogateLossUpdate";   E = "/html/body/div[2]/div/form/div/div[2]/div[1]/div/div[1]/div/div[2]/div/div[1]/div/div/div[3]/div/div/div/div[1]/div/div/div[2]/div[2]/div/table/tbody/tr[2]/td[9]/a" },
    @{ Row = 87; D = "SubrogateLossAmtUpdate"; E = "/html/body/div[2]/div/form/div/div[2]/div[1]/div/div[1]/div/div[2]/div/div[1]/div/div/div[3]/div/div/div/div[4]/div/div/div[2]/div[1]/div[2]/div[1]/div[4]/div/div/input" },
    @{ Row = 88; D = "SubrogateLossAmtSubmit"; E = "/html/body/div[2]/div/form/div/div[2]/div[1]/div/div[1]/div/div[2]/div/div[1]/div/div/div[3]/div/div/div/div[4]/div/div/div[2]/div[2]/button[1]" }
)

foreach ($item in $newRows) {
    $ws.Cells.Item($item.Row, 4).Value = $item.D
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}

# Match the updated view state recorded in the workbook (scroll position and
# active cell selection after the new rows were appended).
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D92").Select()
